$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting the existing rows 44-106 down to 45-107.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the new data record.
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44581
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112022
$ws.Range("G44").Value = "Arveja Verde"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 70
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = 30000
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 1200
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"

# Make sure the date cell keeps the same numeric date format used by the other
# rows in column D (style index 2 in the original workbook).
$ws.Range("D44").NumberFormat = $ws.Range("D45").NumberFormat
